$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update row 2 text (B2): "Invesitgate specular issue" -> "Point lights - specular lighting"
$ws.Range("B2").Value = "Point lights - specular lighting"

# 2. Delete rows 3 and 4 ("Point lights" and "Point lights in model chain"), shifting rows up
$ws.Rows("3:4").Delete()

# 3. Update selection to B3
$ws.Range("B3").Select()
